$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 497.7143
$ws.Range("I12").Value = 297
$ws.Range("J12").Value = 999.5
$ws.Range("K12").Value = 297
$ws.Range("L12").Value = 999.5
$ws.Range("M12").Value = -127
$ws.Range("N12").Value = -1339.5

$ws.Range("H33").Value = 711.82355
$ws.Range("I33").Value = 717.1
$ws.Range("J33").Value = 704.2857
$ws.Range("K33").Value = 717.1
$ws.Range("L33").Value = 704.2857
$ws.Range("M33").Value = -488.1
$ws.Range("N33").Value = -1162.2857

$ws.Range("H62").Value = 2391.6667
$ws.Range("I62").Value = 2391.6667
$ws.Range("K62").Value = 2391.6667
$ws.Range("M62").Value = -1767.6667

$ws.Range("H65").Value = 2391.6667
$ws.Range("I65").Value = 2391.6667
$ws.Range("K65").Value = 11958.3335
$ws.Range("M65").Value = -8838.333500000001

$ws.Range("H112").Value = 3537.1428
$ws.Range("J112").Value = 3752
$ws.Range("L112").Value = 11256
$ws.Range("N112").Value = -13472

$ws.Range("H137").Value = 2320.6924
$ws.Range("I137").Value = 2121.9
$ws.Range("J137").Value = 2983.3333
$ws.Range("K137").Value = 6365.700000000001
$ws.Range("L137").Value = 8949.999899999999
$ws.Range("M137").Value = -3815.700000000001
$ws.Range("N137").Value = -14049.9999

$ws.Range("H138").Value = 2790.9792
$ws.Range("J138").Value = 2666.3948
$ws.Range("L138").Value = 7999.1844
$ws.Range("N138").Value = -18279.1844

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2796.672
$ws.Range("J32").Value = 10000
$ws.Range("L32").Value = 10000
$ws.Range("N32").Value = -10574

$ws.Range("H45").Value = 5550
$ws.Range("I45").Value = 3475
$ws.Range("J45").Value = 6142.857
$ws.Range("K45").Value = 3475
$ws.Range("L45").Value = 6142.857
$ws.Range("M45").Value = -3098
$ws.Range("N45").Value = -6896.857

$ws.Range("H61").Value = 10802
$ws.Range("I61").Value = 9747.056
$ws.Range("K61").Value = 9747.056
$ws.Range("M61").Value = -9535.056

$ws.Range("H74").Value = 5651.9546
$ws.Range("I74").Value = 5306.8423
$ws.Range("K74").Value = 5306.8423
$ws.Range("M74").Value = -4432.8423

$ws.Range("H77").Value = 5651.9546
$ws.Range("I77").Value = 5306.8423
$ws.Range("K77").Value = 26534.2115
$ws.Range("M77").Value = -22166.2115

$ws.Range("H122").Value = 5000
$ws.Range("J122").Value = 5000
$ws.Range("L122").Value = 15000
$ws.Range("N122").Value = -19900

$ws.Range("H136").Value = 10802
$ws.Range("I136").Value = 9747.056
$ws.Range("K136").Value = 29241.168
$ws.Range("M136").Value = -26691.168

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4814.2354
$ws.Range("J31").Value = 5364.778
$ws.Range("L31").Value = 5364.778
$ws.Range("N31").Value = -5954.778

$ws.Range("H34").Value = 4814.2354
$ws.Range("J34").Value = 5364.778
$ws.Range("L34").Value = 5364.778
$ws.Range("N34").Value = -5768.778

$ws.Range("H53").Value = 0
$ws.Range("J53").Value = 0
$ws.Range("L53").ClearContents()
$ws.Range("N53").Value = 0

$ws.Range("H58").Value = 6714.6523
$ws.Range("I58").Value = 4584.5293
$ws.Range("K58").Value = 4584.5293
$ws.Range("M58").Value = -4381.5293

$ws.Range("H107").Value = 741.2308
$ws.Range("I107").Value = 571.8889
$ws.Range("J107").Value = 1122.25
$ws.Range("K107").Value = 571.8889
$ws.Range("L107").Value = 1122.25
$ws.Range("M107").Value = 1348.1111
$ws.Range("N107").Value = -4962.25

$ws.Range("H136").Value = 6714.6523
$ws.Range("I136").Value = 4584.5293
$ws.Range("K136").Value = 13753.5879
$ws.Range("M136").Value = -11203.5879

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 7177107.5
$ws.Range("I4").Value = 13855074
$ws.Range("J4").Value = 4773039.5
$ws.Range("K4").Value = 41565222
$ws.Range("L4").Value = 14319118.5
$ws.Range("M4").Value = -41565110
$ws.Range("N4").Value = -14319342.5

$ws.Range("H12").Value = 855.5
$ws.Range("J12").Value = 1069.125
$ws.Range("L12").Value = 3207.375
$ws.Range("N12").Value = -3553.375

$ws.Range("H29").Value = 558.625
$ws.Range("I29").Value = 499
$ws.Range("J29").Value = 737.5
$ws.Range("K29").Value = 1497
$ws.Range("L29").Value = 2212.5
$ws.Range("M29").Value = -1220
$ws.Range("N29").Value = -2766.5

$ws.Range("H128").Value = 1040756.3
$ws.Range("I128").Value = 1040756.3
$ws.Range("K128").Value = 3122268.9
$ws.Range("M128").Value = -3117288.9

$ws.Range("H132").Value = 1954.45
$ws.Range("J132").Value = 2099.5833
$ws.Range("L132").Value = 18896.2497
$ws.Range("N132").Value = -23956.2497

$ws.Range("H140").Value = 949277.2
$ws.Range("I140").Value = 2444.125
$ws.Range("J140").Value = 3474165.2
$ws.Range("K140").Value = 7332.375
$ws.Range("L140").Value = 10422495.6
$ws.Range("M140").Value = -2152.375
$ws.Range("N140").Value = -10432855.6

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H19").Value = 20152.5
$ws.Range("I19").Value = 0
$ws.Range("K19").Value = 0
$ws.Range("M19").ClearContents()

$ws.Range("H44").Value = 19999.666
$ws.Range("J44").Value = 19999.666
$ws.Range("L44").Value = 19999.666
$ws.Range("N44").Value = -21191.666

$ws.Range("H102").Value = 3840.2354
$ws.Range("I102").Value = 2598.5
$ws.Range("K102").Value = 2598.5
$ws.Range("M102").Value = -976.5

$ws.Range("H113").Value = 4000000
$ws.Range("I113").Value = 4000000
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 4000000
$ws.Range("L113").Value = 0
$ws.Range("M113").ClearContents()
$ws.Range("N113").Value = -3997830

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1678.619
$ws.Range("I16").Value = 1547.2778
$ws.Range("K16").Value = 1547.2778
$ws.Range("M16").Value = -1377.2778

$ws.Range("H22").Value = 1088.909
$ws.Range("I22").Value = 1346.3334
$ws.Range("J22").Value = 780
$ws.Range("K22").Value = 1346.3334
$ws.Range("L22").Value = 780
$ws.Range("M22").Value = -1051.3334
$ws.Range("N22").Value = -1370

$ws.Range("H27").Value = 1088.909
$ws.Range("I27").Value = 1346.3334
$ws.Range("J27").Value = 780
$ws.Range("K27").Value = 1346.3334
$ws.Range("L27").Value = 780
$ws.Range("M27").Value = -1239.3334
$ws.Range("N27").Value = -994

$ws.Range("H40").Value = 3666.4443
$ws.Range("I40").Value = 3474.875
$ws.Range("K40").Value = 3474.875
$ws.Range("M40").Value = -3338.875

$ws.Range("H93").Value = 9680
$ws.Range("I93").Value = 2134.5
$ws.Range("J93").Value = 18483.084
$ws.Range("K93").Value = 2134.5
$ws.Range("L93").Value = 18483.084
$ws.Range("M93").Value = -886.5
$ws.Range("N93").Value = -20979.084

$ws.Range("H136").Value = 4071.795
$ws.Range("I136").Value = 3337.743
$ws.Range("J136").Value = 10494.75
$ws.Range("K136").Value = 10013.229
$ws.Range("L136").Value = 31484.25
$ws.Range("M136").Value = -7463.228999999999
$ws.Range("N136").Value = -36584.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1074.5
$ws.Range("I81").Value = 1113.8572
$ws.Range("J81").Value = 799
$ws.Range("K81").Value = 2227.7144
$ws.Range("L81").Value = 1598
$ws.Range("M81").Value = -1166.7144
$ws.Range("N81").Value = -3720

$ws.Range("H84").Value = 1074.5
$ws.Range("I84").Value = 1113.8572
$ws.Range("J84").Value = 799
$ws.Range("K84").Value = 11138.572
$ws.Range("L84").Value = 7990
$ws.Range("M84").Value = -5834.572
$ws.Range("N84").Value = -18598

$ws.Range("H100").Value = 683.4
$ws.Range("I100").Value = 891
$ws.Range("K100").Value = 1782
$ws.Range("M100").Value = -1241

$ws.Range("H132").Value = 3723.1462
$ws.Range("I132").Value = 3955.8667
$ws.Range("K132").Value = 11867.6001
$ws.Range("M132").Value = -9337.6001
